$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos: value (B/C) changes from the Adriano Siqueira text to the real objectives text ---
$ws.Range("B10").Value = 'Aplicar as técnicas de modelação matemática no estudo de processos de tratamento de águas de abastecimento e residuárias. Fornecer ao aluno condições para uma análise matemática dos sistemas de tratamento de resíduos através de fundamentos de modelagem de fenômenos físicos e bioquímicos. Desenvolver a capacidade de uso de modelos matemáticos na simulação de processos empregados no tratamento de águas.'
$ws.Range("C10").Value = 'Aplicar as técnicas de modelação matemática no estudo de processos de tratamento de águas de abastecimento e residuárias. Fornecer ao aluno condições para uma análise matemática dos sistemas de tratamento de resíduos através de fundamentos de modelagem de fenômenos físicos e bioquímicos. Desenvolver a capacidade de uso de modelos matemáticos na simulação de processos empregados no tratamento de águas.'

# --- Row 13: drop the "Programa resumido:" label from A13 (it moves to A14);
#     B13/C13 get the "229266 - Adriano Francisco Siqueira" value (moved up from row 18) ---
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '229266 - Adriano Francisco Siqueira'
$ws.Range("C13").Value = '229266 - Adriano Francisco Siqueira'
$ws.Rows.Item(13).EntireRow.AutoFit()

# --- Row 14: label becomes "Programa resumido:" (moved up from A13); new B14/C14 hold the summary text ---
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B10:C10").Copy() | Out-Null
$ws.Range("B14:C14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = 'Modelagem matemática de processos de tratamento de águas residuárias: sedimentação, aeração, reatores aeróbios, reatores anaeróbios. Modelos matematicos de processos de tratamento de águas de abastecimento: floculação e filtração. Calibração e validação de modelos.'
$ws.Range("C14").Value = 'Modelagem matemática de processos de tratamento de águas residuárias: sedimentação, aeração, reatores aeróbios, reatores anaeróbios. Modelos matematicos de processos de tratamento de águas de abastecimento: floculação e filtração. Calibração e validação de modelos.'

# --- Row 15: label becomes "Short syllabus:" (moved up from A16); B/C cleared; height 60 ---
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16: label becomes "Programa:"; new B16/C16 hold the full program text ---
$ws.Range("A16").Value = "Programa:"
$ws.Range("B10:C10").Copy() | Out-Null
$ws.Range("B16:C16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = '1- Dinâmica de processos físico-químicos e biológicos. 2- Revisão das equações fundamentais: cinética bio-química e conservação da massa. 3- Fundamentos dos modelos de floculação. 4- Modelos dinâmicos do processo de sedimentação. 5- Fundamentos dos modelos de filtração: a equação de Darcy e os modelos de resistência à filtração. 6- Modelos dinâmicos do processo de oxigenação de águas com e sem consumo simultâneo de oxigênio. 7- Fundamentos do modelo de tratamento de águas residuárias por lodos ativados. 8- Fundamentos dos modelos de digestão anaeróbia. 9- Calibração e validação de modelos.'
$ws.Range("C16").Value = '1- Dinâmica de processos físico-químicos e biológicos. 2- Revisão das equações fundamentais: cinética bio-química e conservação da massa. 3- Fundamentos dos modelos de floculação. 4- Modelos dinâmicos do processo de sedimentação. 5- Fundamentos dos modelos de filtração: a equação de Darcy e os modelos de resistência à filtração. 6- Modelos dinâmicos do processo de oxigenação de águas com e sem consumo simultâneo de oxigênio. 7- Fundamentos do modelo de tratamento de águas residuárias por lodos ativados. 8- Fundamentos dos modelos de digestão anaeróbia. 9- Calibração e validação de modelos.'

# --- Row 17: label becomes "Syllabus:"; height 120 ---
$ws.Range("A17").Value = "Syllabus:"
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18: label becomes "Avaliação:"; B/C cleared ---
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).EntireRow.AutoFit()

# --- Row 19: label becomes "Método:" ---
$ws.Range("A19").Value = "Método:"

# --- Row 20: label becomes "Critério:" ---
$ws.Range("A20").Value = "Critério:"

# --- Row 21: label becomes "Norma de recuperação:"; height 60 ---
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22: label becomes "Bibliografia:"; new B22/C22 hold the bibliography text; height 120 ---
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B10:C10").Copy() | Out-Null
$ws.Range("B22:C22").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").Value = 'Pinto, José Carlos e Lage, Paulo Laranjeira C. Métodos Numéricos em Problemas de Engenharia Química. Rio de Janeiro, E-papers Serviços Editorias, 2001.Weber Jr., W. J. e DiGianno, F.A Process Dynamics in Environmental Systems.New York, J. Wiley & Sons. 1996.Garcia, Claudio. Modelagem e Simulação de Processos Industriais e de Sistemas Eletromecânicos. São Paulo, Edusp. 1997.Dochain, Denis e Vanrolleghem, Peter. A. Dynamical Modelling and Estimation in Wastewater Treatment Processes. London, IWA Publishing, 2001'
$ws.Range("C22").Value = 'Pinto, José Carlos e Lage, Paulo Laranjeira C. Métodos Numéricos em Problemas de Engenharia Química. Rio de Janeiro, E-papers Serviços Editorias, 2001.Weber Jr., W. J. e DiGianno, F.A Process Dynamics in Environmental Systems.New York, J. Wiley & Sons. 1996.Garcia, Claudio. Modelagem e Simulação de Processos Industriais e de Sistemas Eletromecânicos. São Paulo, Edusp. 1997.Dochain, Denis e Vanrolleghem, Peter. A. Dynamical Modelling and Estimation in Wastewater Treatment Processes. London, IWA Publishing, 2001'
$ws.Rows.Item(22).RowHeight = 120

# --- Row 23: new label "Requisitos:"; B/C cleared ---
$ws.Range("A23").Value = "Requisitos:"
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()
$ws.Rows.Item(23).EntireRow.AutoFit()

# --- Row 24: now holds the LOB1006 requirement line (moved up from row 23) ---
$ws.Range("B24").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"

# --- Row 25 (new row): holds the LOT2035 requirement line (moved down from row 24) ---
$ws.Range("B24:C24").Copy() | Out-Null
$ws.Range("B25:C25").PasteSpecial(-4122) | Out-Null
$ws.Range("B25").Value = "LOT2035 -  Tratamento Biológico de Efluentes  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOT2035 -  Tratamento Biológico de Efluentes  (Requisito fraco)`n"
$ws.Rows.Item(25).RowHeight = 30
